$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.594.21"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.463.97"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -1.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("E9").Value = "  +6.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.99%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "2.843.75"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").Value = "2.416.72"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "41.564.70"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.94%  "
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  +4.69%  "
$ws.Range("E27").Value = "  +3.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.54%  "
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("E39").Value = "  +4.98%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "1.960.18"
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("E47").Value = "  +5.58%  "
$ws.Range("D48").Value = "2.702.29"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.171"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "
